$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetTheme = "Pandemia da COVID-19"
$explicacao = "O conteúdo do tema aborda a gestão da Pandemia, a regulação de atividades de linha de frente ou outras atividades da sociedade brasileira que foram afetadas pela COVID-19, como a assistência social em decorrência ao período da crise sanitária e afins."

$lastRow = $ws.UsedRange.Rows.Count

$updated = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $theme = $ws.Cells.Item($r, 4).Value2
    if ($theme -eq $targetTheme) {
        $jCell = $ws.Cells.Item($r, 10)
        $existing = $jCell.Value2
        if (-not $existing) {
            $jCell.Value = $explicacao
            $updated++
        }
    }
}

Write-Host "Updated $updated rows"
